$d = $word.ActiveDocument

$pairs = @(
    @("536÷7=76, 4", "282÷8=35, 2"),
    @("922÷4=230, 2", "495÷9=55, 0"),
    @("428÷3=142, 2", "567÷9=63, 0"),
    @("579÷3=193, 0", "850÷6=141, 4"),
    @("653÷9=72, 5", "995÷5=199, 0"),
    @("102÷9=11, 3", "842÷9=93, 5"),
    @("158÷4=39, 2", "370÷9=41, 1"),
    @("129÷2=64, 1", "894÷8=111, 6"),
    @("775÷7=110, 5", "692÷2=346, 0"),
    @("612÷6=102, 0", "296÷7=42, 2"),
    @("889÷6=148, 1", "253÷4=63, 1"),
    @("189÷2=94, 1", "823÷2=411, 1"),
    @("560÷6=93, 2", "852÷3=284, 0"),
    @("436÷2=218, 0", "893÷8=111, 5"),
    @("559÷6=93, 1", "264÷5=52, 4"),
    @("705÷4=176, 1", "927÷8=115, 7"),
    @("772÷4=193, 0", "562÷6=93, 4"),
    @("646÷4=161, 2", "646÷7=92, 2"),
    @("910÷4=227, 2", "142÷5=28, 2"),
    @("696÷8=87, 0", "756÷9=84, 0"),
    @("583÷9=64, 7", "716÷4=179, 0"),
    @("266÷6=44, 2", "512÷9=56, 8"),
    @("370÷3=123, 1", "228÷6=38, 0"),
    @("928÷4=232, 0", "352÷9=39, 1"),
    @("419÷3=139, 2", "269÷8=33, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
